# Pip_Version_1.docx edit:
#   "...in this folder run "ls -lrt" ... with "chmod"."
# becomes
#   "...in this folder, run "ls -lrt" ... with "chmod.""
#
# i.e. (1) a comma is inserted right after "this folder", and
#      (2) the trailing closing-quote/period pair is swapped from
#          [U+201D "."] to [". " U+201D] (the period moves inside the quote).

$d = $word.ActiveDocument

$lsq = [char]0x201C   # “
$rsq = [char]0x201D   # ”

$target = "script. To list the directories and files in this folder run " + `
          $lsq + "ls -lrt" + $rsq + " and later change the executable permission " + `
          "for the file with " + $lsq + "chmod" + $rsq + "."

$rng = $d.Content
$found = $rng.Find.Execute($target)
if (-not $found) {
    throw "Could not locate the target sentence to edit."
}
$start = $rng.Start
$end = $rng.End

$beforeComma = "script. To list the directories and files in this folder"
$posAfterFolder = $start + $beforeComma.Length

# --- Step 1: insert the comma after "this folder" -----------------------
$commaPoint = $d.Range($posAfterFolder, $posAfterFolder)
$commaPoint.InsertAfter(",")

# Everything from $posAfterFolder onward shifted right by 1 char.
$end = $end + 1

# --- Step 2: swap the trailing ". " -> ". " -------------------------------
$tailStart = $end - 2
$tailRange = $d.Range($tailStart, $end)
if ($tailRange.Text -ne ($rsq + ".")) {
    throw "Unexpected trailing text before edit: [$($tailRange.Text)]"
}
$tailRange.Text = "." + $rsq

# --- Step 3: pin the newly-created segments as their own runs -----------
# (Word/this engine auto-coalesces adjacent runs that share identical
# formatting; toggling a character property on/off forces the run split
# to stick without leaving any lasting formatting change.)
$commaActual = $d.Range($posAfterFolder, $posAfterFolder + 1)
$commaActual.Bold = 1
$commaActual.Bold = 0

$tailActual = $d.Range($tailStart, $tailStart + 2)
$tailActual.Bold = 1
$tailActual.Bold = 0

Write-Output ("Result: " + $d.Range($start, $end).Text)
